$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-28 Thursday", "2025-08-29 Friday"),
    @("87×60=", "51×25="),
    @("76×56=", "67×69="),
    @("47×18=", "76×29="),
    @("16×83=", "28×11="),
    @("61×88=", "78×23="),
    @("47×68=", "71×26="),
    @("59×25=", "85×87="),
    @("52×59=", "27×71="),
    @("83×57=", "90×90="),
    @("28×85=", "65×38="),
    @("12×21=", "68×63="),
    @("59×79=", "93×16="),
    @("60×18=", "89×87="),
    @("60×58=", "90×56="),
    @("41×76=", "31×27="),
    @("63×68=", "57×35="),
    @("53×18=", "74×94="),
    @("19×23=", "15×33="),
    @("38×69=", "33×73="),
    @("75×12=", "69×37="),
    @("47×53=", "84×27="),
    @("91×79=", "24×81="),
    @("37×62=", "96×80="),
    @("74×28=", "85×86="),
    @("80×89=", "99×55=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
